# Fruta / hortaliza, semanal
# Insert two new weekly sample rows (Primera/Segunda) at the top of the
# data block for row 652, pushing the existing history down by 2 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 652:688 down to 654:690, carrying formatting with them.
$ws.Rows("652:653").Insert()

# New row 652 - "Primera" quality sample for the new week.
$ws.Cells.Item(652, 1).Value = 3
$ws.Cells.Item(652, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(652, 3).Value = "Coquimbo"
$ws.Cells.Item(652, 4).Value = 44706
$ws.Cells.Item(652, 5).Value = 5
$ws.Cells.Item(652, 6).Value = 100114014
$ws.Cells.Item(652, 7).Value = "Betarraga"
$ws.Cells.Item(652, 8).Value = "Sin especificar"
$ws.Cells.Item(652, 9).Value = "Primera"
$ws.Cells.Item(652, 10).Value = 2300
$ws.Cells.Item(652, 11).Value = 650
$ws.Cells.Item(652, 12).Value = 700
$ws.Cells.Item(652, 13).Value = 676
$ws.Cells.Item(652, 14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(652, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(652, 16).Value = 169
$ws.Cells.Item(652, 17).Value = 4
$ws.Cells.Item(652, 18).Value = "Hortaliza"

# New row 653 - "Segunda" quality sample for the new week.
$ws.Cells.Item(653, 1).Value = 3
$ws.Cells.Item(653, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(653, 3).Value = "Coquimbo"
$ws.Cells.Item(653, 4).Value = 44706
$ws.Cells.Item(653, 5).Value = 5
$ws.Cells.Item(653, 6).Value = 100114014
$ws.Cells.Item(653, 7).Value = "Betarraga"
$ws.Cells.Item(653, 8).Value = "Sin especificar"
$ws.Cells.Item(653, 9).Value = "Segunda"
$ws.Cells.Item(653, 10).Value = 1150
$ws.Cells.Item(653, 11).Value = 450
$ws.Cells.Item(653, 12).Value = 450
$ws.Cells.Item(653, 13).Value = 450
$ws.Cells.Item(653, 14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(653, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(653, 16).Value = 112
$ws.Cells.Item(653, 17).Value = 4
$ws.Cells.Item(653, 18).Value = "Hortaliza"
